$d = $word.ActiveDocument

# Paragraph 2 holds the m2doc field ( m:    comment ... ) together with the
# two generated "hint" runs (the blue "You might want..." one and the red
# "Couldn't find..." one). The parser used to build this as a real Word
# field (fldChar begin / instrText* / fldChar end); the new
# TokenIteratorFieldRewriterSplit instead emits the token as plain literal
# text runs ("{", "m", ":", spaces, the comment text, "}") while keeping the
# bookmark that sat in the middle, and moves the red "Couldn't find..." hint
# after the blue one.

$p2 = $d.Paragraphs(2)

# Remove the whole field (fldChar begin/end + all instrText runs). This also
# collaterally removes the bookmark and the red hint runs that were nested
# between the field's begin/end markers - we recreate all of that below from
# an explicit OOXML fragment so the exact target markup is produced.
$f = $d.Fields.Item(1)
$f.Delete()

# Build the exact replacement markup for the whole (now field-less)
# paragraph 2 and drop it in via InsertXML, which replaces the contents of
# the target range with precisely the runs/bookmarks we specify (no
# unwanted run-merging or implicit formatting inheritance).
$replacementXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve">   </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>comment some important comment</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:color w:val="0000FF"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r><w:r><w:rPr><w:color w:val="0000FF"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>You might want to replace m:   comment by m:comment</w:t></w:r><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>Couldn't find the 'comment' variable</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$target.InsertXML($replacementXml)
